# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Only column G (header "K") values change for rows 2-17; all other data is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 6
    3  = 9
    4  = 4
    5  = 4
    6  = 3
    7  = 3
    8  = 6
    9  = 6
    10 = 4
    11 = 5
    12 = 3
    13 = 8
    14 = 3
    15 = 8
    16 = 5
    17 = 4
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
